$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计") ---
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Copy the header row + column-A formatting from the "2021-Q4" sheet, which
# already has the bold/centered/bordered style used across these fund sheets.
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2:A7").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$fundRows = @(
    @(0, "000727", "融通健康产业灵活配置混合A", "15.30", "94.68", "8.08", "1.2362", 4),
    @(1, "009805", "国泰医药健康股票A", "12.14", "92.49", "4.34", "0.5269", 10),
    @(2, "009274", "融通健康产业灵活配置混合C", "3.16", "94.68", "8.08", "0.2553", 4),
    @(3, "011326", "国泰医药健康股票C", "1.09", "92.49", "4.34", "0.0473", 10),
    @(4, "011807", "平安研究精选混合型证券投资基金A", "1.19", "93.59", "2.74", "0.0326", 10),
    @(5, "011808", "平安研究精选混合型证券投资基金C", "0.49", "93.59", "2.74", "0.0134", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 2).Style = "Normal"

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 4).Style = "Normal"

    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 5).Style = "Normal"

    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 6).Style = "Normal"

    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 7).Style = "Normal"

    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# --- 2. Update "总计" sheet: insert a new first data row for 2022-Q1 ---
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# Restore the formatting that "Insert" drops on the new blank row by
# cloning it from the (now pushed-down) row beneath it.
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 6
$totalSheet.Cells.Item(2, 4).Value = 2.11

# Renumber the A column index (0-based order) for the rows that got shifted down
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
